$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell -> new text value. Values must remain plain text (matches the
# original inlineStr cells), so we force Text format before assigning
# and restore the default "Normal" style afterwards so cell styling is
# untouched (no stray quotePrefix / custom numFmt left behind).
$updates = @{
    'D2' = '303.08'
    'E2' = '2.21%'
    'D3' = '31.73'
    'E3' = '0.66%'
    'D4' = '5.168'
    'E4' = '1.07%'
    'D5' = '0.07818'
    'E5' = '4.11%'
    'D6' = '2.430'
    'E6' = '42.74%'
    'D7' = '7.949'
    'E7' = '2.74%'
    'D8' = '3.869'
    'E8' = '1.69%'
    'D9' = '0.9104'
    'E9' = '-2.23%'
    'D10' = '0.1728'
    'E10' = '2.24%'
    'D11' = '0.07332'
    'E11' = '-0.36%'
    'D12' = '0.08141'
    'E12' = '2.73%'
    'D13' = '0.03035'
    'D14' = '0.09939'
    'E14' = '0.47%'
    'D15' = '0.001505'
    'E15' = '0.62%'
    'D16' = '0.006093'
    'E16' = '-5.78%'
    'D17' = '3.498'
    'E17' = '1.54%'
    'D18' = '2.244'
    'E18' = '0.75%'
    'E20' = '0.71%'
    'D21' = '4.673'
    'E21' = '2.31%'
    'D22' = '0.04653'
    'E22' = '-0.20%'
    'E23' = '0.37%'
    'D24' = '0.001260'
    'E24' = '3.33%'
    'D25' = '0.004520'
    'E25' = '2.26%'
    'E26' = '3.69%'
    'D27' = '0.0002740'
    'E27' = '45.96%'
    'D39' = '0.01783'
    'E39' = '6.75%'
    'D40' = '0.04554'
    'E40' = '2.18%'
    'D41' = '0.007272'
    'E41' = '3.76%'
    'D42' = '0.1363'
    'E42' = '2.76%'
    'D43' = '0.002239'
    'E43' = '8.57%'
    'E44' = '-4.25%'
    'D45' = '0.00006476'
    'E45' = '7.83%'
    'D46' = '0.00000000749'
    'E46' = '-0.08%'
    'E47' = '-57.48%'
    'D49' = '0.00002099'
    'E49' = '-0.08%'
    'D50' = '0.0001999'
    'E50' = '-0.01%'
}

foreach ($cell in $updates.Keys) {
    $rng = $ws.Range($cell)
    $rng.NumberFormat = "@"
    $rng.Value = $updates[$cell]
    $rng.Style = "Normal"
}

